# Adapt column header formatting to respective input file names:
#   *_old -> *_FV2410   (left-hand "before" block of the AHB diff)
#   *_new -> *_FV2504   (right-hand "after" block of the AHB diff)
# Then turn the data range into a proper Excel Table and freeze the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$baseNames = @(
    "Segmentname",
    "Segmentgruppe",
    "Segment",
    "Datenelement",
    "Segment ID",
    "Code",
    "Qualifier",
    "Beschreibung",
    "Bedingungsausdruck",
    "Bedingung"
)

# Columns A-J (1-10): "<name>_FV2410"
for ($i = 0; $i -lt $baseNames.Count; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = "$($baseNames[$i])_FV2410"
}

# Column K (11): unchanged "diff" marker column
$ws.Cells.Item(1, 11).Value = "diff"

# Columns L-U (12-21): "<name>_FV2504"
for ($i = 0; $i -lt $baseNames.Count; $i++) {
    $ws.Cells.Item(1, $i + 12).Value = "$($baseNames[$i])_FV2504"
}

# Turn the used range into an Excel Table (adds xl/tables/table1.xml + tableParts)
$dataRange = $ws.Range("A1:U57")
$table = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $dataRange, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$table.Name = "Table1"
$table.TableStyle = ""

# Freeze the header row (pane split after row 1)
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
